# Minor edit: Exporting tables for appendix
# Updates the data table on Sheet1 with refreshed computation results and
# appends a new "Italy" row (existing row 8 "Spain" data moves to row 9,
# old row 8 is replaced with new "USA_WA" data, and a brand new "Italy"
# row 10 is appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style used by the existing Date column (column B) cells so the
# newly appended row keeps the same yyyy-mm-dd number format.
$dateNumberFormat = $ws.Range("B2").NumberFormat

# --- Row 2: Germany ---
$ws.Range("B2").Value = 43919
$ws.Range("C2").Value = 0.006754162448951098
$ws.Range("D2").Value = 0.009464650821448846
$ws.Range("E2").Value = -0.0004602093717535852
$ws.Range("F2").Value = 0.009924860193202432
$ws.Range("G2").Value = 0.0443145198859854
$ws.Range("H2").Value = 0.9556854801140146

# --- Row 3: USA_NYC ---
$ws.Range("D3").Value = 0.003947726222943577
$ws.Range("E3").Value = -0.0001608464805649595
$ws.Range("F3").Value = 0.004108572703508538
$ws.Range("G3").Value = 0.03767408952603575
$ws.Range("H3").Value = 0.9623259104739642

# --- Row 4: SouthKorea ---
$ws.Range("B4").Value = 43920
$ws.Range("C4").Value = 0.01621881327039994

# --- Row 5: US ---
$ws.Range("D5").Value = -0.001104083507632008
$ws.Range("E5").Value = -0.01130159718274817
$ws.Range("F5").Value = 0.01019751367511616
$ws.Range("G5").Value = 0.5256774225439219
$ws.Range("H5").Value = 0.474322577456078

# --- Row 6: China ---
$ws.Range("D6").Value = -0.006546490597355136
$ws.Range("E6").Value = -0.0009030758576810743
$ws.Range("F6").Value = -0.005643414739674061
$ws.Range("G6").Value = 0.1379480874907127
$ws.Range("H6").Value = 0.8620519125092873

# --- Row 7: France ---
$ws.Range("B7").Value = 43914
$ws.Range("C7").Value = 0.03983587515221891
$ws.Range("D7").Value = -0.02361706188181897
$ws.Range("E7").Value = -0.02098271874324288
$ws.Range("F7").Value = -0.002634343138576091
$ws.Range("G7").Value = 0.8884559327591854
$ws.Range("H7").Value = 0.1115440672408145

# --- Row 8: was Spain, now USA_WA ---
$ws.Range("A8").Value = "USA_WA"
$ws.Range("B8").Value = 43919
$ws.Range("C8").Value = 0.04349561395511622
$ws.Range("D8").Value = -0.02727680068471627
$ws.Range("E8").Value = -0.01993017427657376
$ws.Range("F8").Value = -0.007346626408142511
$ws.Range("G8").Value = 0.7306639259838501
$ws.Range("H8").Value = 0.26933607401615

# --- Row 9: new Spain data ---
$ws.Range("A9").Value = "Spain"
$ws.Range("B9").Value = 43919
$ws.Range("C9").Value = 0.08615529080345091
$ws.Range("D9").Value = -0.06993647753305096
$ws.Range("E9").Value = -0.03964027430057621
$ws.Range("F9").Value = -0.03029620323247475
$ws.Range("G9").Value = 0.5668039869729326
$ws.Range("H9").Value = 0.4331960130270674

# --- Row 10: new Italy row (appended) ---
$ws.Range("A10").Value = "Italy"
$ws.Range("B10").Value = 43919
$ws.Range("B10").NumberFormat = $dateNumberFormat
$ws.Range("C10").Value = 0.1063320056230004
$ws.Range("D10").Value = -0.09011319235260046
$ws.Range("E10").Value = -0.0518290945368735
$ws.Range("F10").Value = -0.03828409781572694
$ws.Range("G10").Value = 0.5751554593035994
$ws.Range("H10").Value = 0.4248445406964007
